$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.302.75'
$ws.Range('E2').Value = '  -3.30%  '
$ws.Range('D3').Value = '1.930.12'
$ws.Range('E3').Value = '  -3.79%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '248.15'
$ws.Range('E5').Value = '  -4.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.7279'
$ws.Range('E6').Value = '  -7.76%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3321'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '28.14'
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06927'
$ws.Range('E10').Value = '  -2.40%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8041'
$ws.Range('E11').Value = '  -6.47%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08045'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').Value = '1.932.12'
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.417'
$ws.Range('E14').Value = '  -3.58%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '94.73'
$ws.Range('E15').Value = '  -6.62%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.58'
$ws.Range('E16').Value = '  -2.47%  '
$ws.Range('D17').Value = '30.308.45'
$ws.Range('E17').Value = '  -3.32%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '253.15'
$ws.Range('E18').Value = '  -8.22%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000008215'
$ws.Range('E19').Value = '  +2.69%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.804'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').Value = '2.188.89'
$ws.Range('E21').Value = '  -3.54%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9994'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -4.62%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.737'
$ws.Range('E25').Value = '  -4.13%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '159.22'
$ws.Range('E26').Value = '  -3.55%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.463'
$ws.Range('E27').Value = '  +2.24%  '
$ws.Range('E28').Value = '  -4.53%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.1339'
$ws.Range('E29').Value = '  -11.82%  '
$ws.Range('E30').Value = '  -4.31%  '
$ws.Range('E31').Value = '  -1.74%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.412'
$ws.Range('E32').Value = '  -4.83%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.197'
$ws.Range('E33').Value = '  -5.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05129'
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.221'
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7440'
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.752'
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01983'
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.834'
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.641'
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '79.01'
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4472'
$ws.Range('E42').Value = '  -6.27%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.997'
$ws.Range('E43').Value = '  -7.58%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9998'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8364'
$ws.Range('E45').Value = '  -2.66%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '102.22'
$ws.Range('E46').Value = '  -4.97%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.794'
$ws.Range('E47').Value = '  -2.55%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.314'
$ws.Range('E48').Value = '  -6.25%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '36.62'
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05959'
$ws.Range('E51').Value = '  -0.47%  '
